$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-21 Thursday" "2024-11-22 Friday"

Replace-Text "90×42=" "97×74="
Replace-Text "95×73=" "96×70="
Replace-Text "93×52=" "76×36="
Replace-Text "55×86=" "43×51="
Replace-Text "94×45=" "36×85="

Replace-Text "59×85=" "59×41="
Replace-Text "64×60=" "81×42="
Replace-Text "42×35=" "32×70="
Replace-Text "38×39=" "28×61="
Replace-Text "91×15=" "32×71="

Replace-Text "59×86=" "32×96="
Replace-Text "82×38=" "84×87="
Replace-Text "70×49=" "92×84="
Replace-Text "33×71=" "81×80="
Replace-Text "14×98=" "71×80="

Replace-Text "70×40=" "86×31="
Replace-Text "64×27=" "28×80="
Replace-Text "35×72=" "99×53="
Replace-Text "19×65=" "33×32="
Replace-Text "12×56=" "61×93="

Replace-Text "79×29=" "86×32="
Replace-Text "21×78=" "31×91="
Replace-Text "43×85=" "96×91="
Replace-Text "42×31=" "75×15="
Replace-Text "61×20=" "31×88="

Write-Output "Done"
